$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Long Cypher query text blocks (verbatim, single-quoted here-strings) ----
$statQueryText = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (f:file)-[*]->(samp:sample)-->(c)
MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp,demo, c, s, p, diag
WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1', 'T2N1M0', 'T2N1M1', 'T3N0M0', 'T3N1M0', 'T3N1M1'] 
OPTIONAL MATCH (samp:sample)-->(c)
RETURN
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

$casesQueryText = @'
 MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis) 
 MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1', 'T2N1M0', 'T2N1M1', 'T3N0M0', 'T3N1M0', 'T3N1M1'] 
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co, demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  
       coalesce(c.case_id, '') AS `Case ID`,
       coalesce(s.clinical_study_designation, '') AS `Study Code`,
       coalesce(s.clinical_study_type, '') AS  `Study Type`,
       coalesce(demo.breed, '') AS Breed ,
       coalesce(diag.disease_term, '') AS Diagnosis ,
       coalesce(diag.stage_of_disease, '') AS `Stage of Disease`,
       CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END AS Age,
       coalesce(demo.sex, '') AS Sex,
       coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
       coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
       coalesce(diag.best_response, '') AS `Response to Treatment`,
       coalesce(co.cohort_description, '') AS `Cohort`
Order by c.case_id LIMIT 100        
'@

$samplesQueryText = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (samp:sample)-->(c)<--(diag:diagnosis) 
WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1', 'T2N1M0', 'T2N1M1', 'T3N0M0', 'T3N1M0', 'T3N1M1'] 
 WITH DISTINCT samp AS samp, c, demo, diag
RETURN  coalesce(samp.sample_id, '') AS `Sample ID`, 
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(samp.sample_site, '') AS `Sample Site`,
        coalesce(samp.summarized_sample_type, '') AS `Sample Type`,
        coalesce(samp.specific_sample_pathology, '') AS `Pathology/Morphology`,
        coalesce(samp.tumor_grade, '') AS `Tumor Grade`,
        coalesce(samp.sample_chronology, '') AS `Sample Chronology`,
        coalesce(samp.percentage_tumor, '') AS `Percentage Tumor`,
        coalesce(samp.necropsy_sample, '') AS `Necropsy Sample`,
        coalesce(samp.sample_preservation, '') AS `Sample Preservation`
'@

$filesQueryText = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1', 'T2N1M0', 'T2N1M1', 'T3N0M0', 'T3N1M0', 'T3N1M1'] 
WITH DISTINCT f, parent, c, demo, diag, s
OPTIONAL MATCH (f)-[*]->(samp:sample)
OPTIONAL MATCH (s:study)<--(c)<--(diag:diagnosis)<-[*]-(samp)
WITH
        f, parent, c, demo, diag, s, samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, parent, c, demo, diag, s, samp,
        f.file_size /(1024^i) AS value, 
        10^precision AS factor,
        units[i] as unit
WITH    
        f, parent, c, demo, diag, s, samp, unit,
        round(factor * value)/factor AS size
RETURN coalesce(f.file_name, '') AS `File Name`, 
 coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`, 
      CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
   coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis 
        Order By f.file_name LIMIT 100
'@

$studyFilesQueryText = @'
MATCH (f:file)-->(s:study)
MATCH (s)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (sf:file)-->(s)
MATCH (s)<--(c)
MATCH (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['UBC02'] and diag.stage_of_disease in [ 'T2N0M0', 'T2N0M1', 'T2N1M0', 'T2N1M1', 'T3N0M0', 'T3N1M0', 'T3N1M1'] 
WITH DISTINCT f,  s, c
WITH
        f, c,  s,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH
        f, c,  s,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH
        f,  c,   s, unit,
        round(factor * value)/factor AS size
RETURN DISTINCT
  coalesce(f.file_name, '') AS `File Name`,
  coalesce(f.file_type, '') AS `File Type`,
  coalesce("study", '') AS `Association`,
  coalesce(f.file_description, '') AS `Description`,
  coalesce(f.file_format, '') AS  Format,
  CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
  coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# ---- Row 2: CasesTab query updated, stats query updated ----
$ws.Range("B2").Value = $casesQueryText
$ws.Range("C2").Value = $statQueryText
$ws.Rows.Item(2).RowHeight = 345

# ---- Row 3: SamplesTab query updated, stats query updated ----
$ws.Range("B3").Value = $samplesQueryText
$ws.Range("C3").Value = $statQueryText
$ws.Rows.Item(3).RowHeight = 255

# ---- Row 4: FilesTab query updated, stats query updated ----
$ws.Range("B4").Value = $filesQueryText
$ws.Range("C4").Value = $statQueryText
$ws.Rows.Item(4).RowHeight = 409.5

# ---- Row 5 (new): StudyFilesTab ----
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $studyFilesQueryText
$ws.Range("C5").Value = $statQueryText
$ws.Range("D5").Value = "TC22_Canine_StudyUBC02-AllBreeds_StageOfDisease_Neo4jData.xlsx"
$ws.Range("E5").Value = "TC22_Canine_StudyUBC02-AllBreeds_StageOfDisease_WebData.xlsx"
$ws.Range("B5:C5").WrapText = $true
$ws.Rows.Item(5).RowHeight = 409.5

# ---- Update sheet view to reflect the newly active/selected row ----
[void]$ws.Range("B5").Select()
